$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows above the current row 82, shifting existing rows 82:104 down to 88:110.
$ws.Range("A82:T87").Insert(4)

# Constant columns shared by every Damasco / Mercado Mayorista Lo Valledor de Santiago record.
$mercadoId = 6
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$tipo = "Fruta"
$productoId = 100103
$producto = "Frutos de hueso (carozo)"
$categoriaId = 100103003
$categoria = "Damasco"

# New weekly rows (fecha serial 44551) inserted at rows 82-87.
$rows = @(
    @{ Row=82; Fecha=44551; Variedad="Dina";    Calidad="Especial"; Volumen=50;  Min=16000;  Max=16000;  Prom=16000;  Unidad="$/caja 16 kilos";    Origen="Región Metropolitana"; PrecioKg=1000; KgUnidad=16 }
    @{ Row=83; Fecha=44551; Variedad="Dina";    Calidad="Primera";  Volumen=85;  Min=14000;  Max=14000;  Prom=14000;  Unidad="$/caja 16 kilos";    Origen="Región Metropolitana"; PrecioKg=875;  KgUnidad=16 }
    @{ Row=84; Fecha=44551; Variedad="Dina";    Calidad="Segunda";  Volumen=80;  Min=11000;  Max=11000;  Prom=11000;  Unidad="$/caja 16 kilos";    Origen="Región Metropolitana"; PrecioKg=688;  KgUnidad=16 }
    @{ Row=85; Fecha=44551; Variedad="Modesto"; Calidad="Especial"; Volumen=8;   Min=400000; Max=400000; Prom=400000; Unidad="$/bins (500 kilos)"; Origen="Región Metropolitana"; PrecioKg=800;  KgUnidad=500 }
    @{ Row=86; Fecha=44551; Variedad="Modesto"; Calidad="Primera";  Volumen=31;  Min=330000; Max=350000; Prom=340323; Unidad="$/bins (500 kilos)"; Origen="Región Metropolitana"; PrecioKg=681;  KgUnidad=500 }
    @{ Row=87; Fecha=44551; Variedad="Modesto"; Calidad="Segunda";  Volumen=12;  Min=270000; Max=270000; Prom=270000; Unidad="$/bins (500 kilos)"; Origen="Región Metropolitana"; PrecioKg=540;  KgUnidad=500 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Cells.Item($i, 1).Value = $mercadoId
    $ws.Cells.Item($i, 2).Value = $mercado
    $ws.Cells.Item($i, 3).Value = $region
    $ws.Cells.Item($i, 4).Value = $r.Fecha
    $ws.Cells.Item($i, 5).Value = $codreg
    $ws.Cells.Item($i, 6).Value = $tipo
    $ws.Cells.Item($i, 7).Value = $productoId
    $ws.Cells.Item($i, 8).Value = $producto
    $ws.Cells.Item($i, 9).Value = $categoriaId
    $ws.Cells.Item($i, 10).Value = $categoria
    $ws.Cells.Item($i, 11).Value = $r.Variedad
    $ws.Cells.Item($i, 12).Value = $r.Calidad
    $ws.Cells.Item($i, 13).Value = $r.Volumen
    $ws.Cells.Item($i, 14).Value = $r.Min
    $ws.Cells.Item($i, 15).Value = $r.Max
    $ws.Cells.Item($i, 16).Value = $r.Prom
    $ws.Cells.Item($i, 17).Value = $r.Unidad
    $ws.Cells.Item($i, 18).Value = $r.Origen
    $ws.Cells.Item($i, 19).Value = $r.PrecioKg
    $ws.Cells.Item($i, 20).Value = $r.KgUnidad
}
